$wb = $excel.ActiveWorkbook
$w = $wb.Windows.Item(1)
$w.Left = 10220
$w.Top = 6420
$w.Width = 34140
$w.Height = 16460
Write-Output "done"
